# ============================================================
# Weekly CompStat (33rd Precinct) refresh: new reporting week,
# updated volume number, and refreshed crime-complaint figures.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 32   Number  18" -> "...19" ---
$hdr = $ws.Cells.Item(8,1)
$hdr.Characters(21,2).Text = "19"

# --- Header: report week dates ---
$wk = $ws.Cells.Item(9,3)
$wk.Characters(47,8).Text = "5/11/2025"
$wk.Characters(27,9).Text = "5/5/2025"

# --- Column E got a hair wider to fit the new figures ---
$ws.Columns.Item(5).ColumnWidth = 6.71

# --- Cells whose content switched between a number and the
#     literal placeholder text ('0' / '***.*'): copy format+type
#     from a same-shaped anchor cell in the (untouched) row above,
#     then set the real value for the numeric ones. ---
$ws.Cells.Item(14,7).Copy($ws.Cells.Item(15,4))
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(14,8).Copy($ws.Cells.Item(15,5))
$ws.Cells.Item(15,5).Value = -100
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(15,6))
$ws.Cells.Item(14,7).Copy($ws.Cells.Item(15,7))
$ws.Cells.Item(15,7).Value = 1
$ws.Cells.Item(14,8).Copy($ws.Cells.Item(15,8))
$ws.Cells.Item(15,8).Value = -100
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(18,4))
$ws.Cells.Item(14,5).Copy($ws.Cells.Item(18,5))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(22,3))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(22,4))
$ws.Cells.Item(14,5).Copy($ws.Cells.Item(22,5))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(23,3))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(23,4))
$ws.Cells.Item(14,5).Copy($ws.Cells.Item(23,5))
$ws.Cells.Item(14,7).Copy($ws.Cells.Item(27,4))
$ws.Cells.Item(27,4).Value = 2
$ws.Cells.Item(14,8).Copy($ws.Cells.Item(27,5))
$ws.Cells.Item(27,5).Value = -100
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(27,6))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(28,3))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(28,4))
$ws.Cells.Item(14,5).Copy($ws.Cells.Item(28,5))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(31,3))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(31,6))
$ws.Cells.Item(14,3).Copy($ws.Cells.Item(31,9))

# --- Remaining cells: value-only refresh (style untouched) ---
$ws.Cells.Item(15,10).Value = 4
$ws.Cells.Item(15,11).Value = 50
$ws.Cells.Item(15,13).Value = 0
$ws.Cells.Item(16,3).Value = 2
$ws.Cells.Item(16,5).Value = -60
$ws.Cells.Item(16,6).Value = 12
$ws.Cells.Item(16,7).Value = 18
$ws.Cells.Item(16,8).Value = -33.333333333333
$ws.Cells.Item(16,9).Value = 36
$ws.Cells.Item(16,10).Value = 63
$ws.Cells.Item(16,11).Value = -42.857142857142
$ws.Cells.Item(16,12).Value = -21.739130434782
$ws.Cells.Item(16,13).Value = -53.846153846153
$ws.Cells.Item(17,3).Value = 11
$ws.Cells.Item(17,4).Value = 4
$ws.Cells.Item(17,5).Value = 175
$ws.Cells.Item(17,6).Value = 26
$ws.Cells.Item(17,7).Value = 17
$ws.Cells.Item(17,8).Value = 52.941176470588
$ws.Cells.Item(17,9).Value = 87
$ws.Cells.Item(17,10).Value = 88
$ws.Cells.Item(17,11).Value = -1.136363636363
$ws.Cells.Item(17,12).Value = -5.434782608695
$ws.Cells.Item(17,13).Value = 40.322580645161
$ws.Cells.Item(18,7).Value = 7
$ws.Cells.Item(18,8).Value = -42.857142857142
$ws.Cells.Item(18,12).Value = -43.396226415094
$ws.Cells.Item(18,13).Value = -16.666666666666
$ws.Cells.Item(19,3).Value = 5
$ws.Cells.Item(19,4).Value = 2
$ws.Cells.Item(19,5).Value = 150
$ws.Cells.Item(19,6).Value = 33
$ws.Cells.Item(19,7).Value = 28
$ws.Cells.Item(19,8).Value = 17.857142857142
$ws.Cells.Item(19,9).Value = 116
$ws.Cells.Item(19,10).Value = 133
$ws.Cells.Item(19,11).Value = -12.781954887218
$ws.Cells.Item(19,12).Value = 14.851485148514
$ws.Cells.Item(19,13).Value = 23.404255319148
$ws.Cells.Item(20,3).Value = 3
$ws.Cells.Item(20,4).Value = 4
$ws.Cells.Item(20,5).Value = -25
$ws.Cells.Item(20,6).Value = 13
$ws.Cells.Item(20,7).Value = 7
$ws.Cells.Item(20,8).Value = 85.714285714285
$ws.Cells.Item(20,9).Value = 47
$ws.Cells.Item(20,10).Value = 34
$ws.Cells.Item(20,11).Value = 38.235294117647
$ws.Cells.Item(20,12).Value = 6.818181818181
$ws.Cells.Item(20,13).Value = 74.074074074074
$ws.Cells.Item(21,3).Value = 22
$ws.Cells.Item(21,4).Value = 16
$ws.Cells.Item(21,5).Value = 37.5
$ws.Cells.Item(21,6).Value = 88
$ws.Cells.Item(21,7).Value = 79
$ws.Cells.Item(21,8).Value = 11.392405063291
$ws.Cells.Item(21,9).Value = 322
$ws.Cells.Item(21,10).Value = 366
$ws.Cells.Item(21,11).Value = -12.021857923497
$ws.Cells.Item(21,12).Value = -5.294117647058
$ws.Cells.Item(21,13).Value = 5.573770491803
$ws.Cells.Item(22,13).Value = 0
$ws.Cells.Item(23,6).Value = 2
$ws.Cells.Item(23,8).Value = 100
$ws.Cells.Item(24,3).Value = 14
$ws.Cells.Item(24,4).Value = 17
$ws.Cells.Item(24,5).Value = -17.647058823529
$ws.Cells.Item(24,6).Value = 48
$ws.Cells.Item(24,7).Value = 89
$ws.Cells.Item(24,8).Value = -46.067415730337
$ws.Cells.Item(24,9).Value = 294
$ws.Cells.Item(24,10).Value = 316
$ws.Cells.Item(24,11).Value = -6.962025316455
$ws.Cells.Item(24,12).Value = -13.529411764705
$ws.Cells.Item(24,13).Value = 65.168539325842
$ws.Cells.Item(25,3).Value = 3
$ws.Cells.Item(25,4).Value = 5
$ws.Cells.Item(25,5).Value = -40
$ws.Cells.Item(25,6).Value = 8
$ws.Cells.Item(25,7).Value = 25
$ws.Cells.Item(25,8).Value = -68
$ws.Cells.Item(25,9).Value = 74
$ws.Cells.Item(25,10).Value = 74
$ws.Cells.Item(25,11).Value = 0
$ws.Cells.Item(25,12).Value = -17.777777777777
$ws.Cells.Item(26,3).Value = 10
$ws.Cells.Item(26,4).Value = 6
$ws.Cells.Item(26,5).Value = 66.666666666666
$ws.Cells.Item(26,6).Value = 44
$ws.Cells.Item(26,8).Value = 69.230769230769
$ws.Cells.Item(26,9).Value = 141
$ws.Cells.Item(26,10).Value = 117
$ws.Cells.Item(26,11).Value = 20.512820512820
$ws.Cells.Item(26,12).Value = -5.369127516778
$ws.Cells.Item(26,13).Value = -3.424657534246
$ws.Cells.Item(27,7).Value = 3
$ws.Cells.Item(27,8).Value = -100
$ws.Cells.Item(27,10).Value = 7
$ws.Cells.Item(27,11).Value = -14.285714285714
$ws.Cells.Item(28,6).Value = 2
$ws.Cells.Item(28,7).Value = 4
$ws.Cells.Item(28,8).Value = -50
$ws.Cells.Item(31,11).Value = -100
$ws.Cells.Item(31,12).Value = -100
